$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1374.125
$ws.Range("I41").Value = 1300
$ws.Range("J41").Value = 1418.6
$ws.Range("K41").Value = 1300
$ws.Range("L41").Value = 1418.6
$ws.Range("M41").Value = -860
$ws.Range("N41").Value = -2298.6
$ws.Range("H53").Value = 13859.533
$ws.Range("I53").Value = 570.2857
$ws.Range("K53").Value = 570.2857
$ws.Range("M53").Value = 66.71429999999998
$ws.Range("H86").Value = 10087.714
$ws.Range("I86").Value = 9528.5
$ws.Range("J86").Value = 10833.333
$ws.Range("K86").Value = 9528.5
$ws.Range("L86").Value = 10833.333
$ws.Range("M86").Value = -8405.5
$ws.Range("N86").Value = -13079.333
$ws.Range("H89").Value = 10087.714
$ws.Range("I89").Value = 9528.5
$ws.Range("J89").Value = 10833.333
$ws.Range("K89").Value = 47642.5
$ws.Range("L89").Value = 54166.665
$ws.Range("M89").Value = -42026.5
$ws.Range("N89").Value = -65398.665
$ws.Range("H107").Value = 27861474
$ws.Range("I107").Value = 47619384
$ws.Range("K107").Value = 47619384
$ws.Range("M107").Value = -47617464
$ws.Range("H137").Value = 50432.137
$ws.Range("I137").Value = 70301.08
$ws.Range("K137").Value = 210903.24
$ws.Range("M137").Value = -208353.24
$ws.Range("H138").Value = 3287.123
$ws.Range("J138").Value = 3532.04
$ws.Range("L138").Value = 10596.12
$ws.Range("N138").Value = -20876.12

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5520.44
$ws.Range("I61").Value = 6308.121
$ws.Range("J61").Value = 3991.4119
$ws.Range("K61").Value = 6308.121
$ws.Range("L61").Value = 3991.4119
$ws.Range("M61").Value = -6096.121
$ws.Range("N61").Value = -4415.4119
$ws.Range("H74").Value = 13913.341
$ws.Range("I74").Value = 3749.5754
$ws.Range("K74").Value = 3749.5754
$ws.Range("M74").Value = -2875.5754
$ws.Range("H77").Value = 13913.341
$ws.Range("I77").Value = 3749.5754
$ws.Range("K77").Value = 18747.877
$ws.Range("M77").Value = -14379.877
$ws.Range("H97").Value = 1766837.4
$ws.Range("J97").Value = 2085
$ws.Range("L97").Value = 2085
$ws.Range("N97").Value = -3077
$ws.Range("H135").Value = 333355680
$ws.Range("J135").Value = 333355680
$ws.Range("L135").Value = 333355680
$ws.Range("N135").Value = -333365820
$ws.Range("H136").Value = 5520.44
$ws.Range("I136").Value = 6308.121
$ws.Range("J136").Value = 3991.4119
$ws.Range("K136").Value = 18924.363
$ws.Range("L136").Value = 11974.2357
$ws.Range("M136").Value = -16374.363
$ws.Range("N136").Value = -17074.2357
$ws.Range("H138").Value = 108497.75
$ws.Range("J138").Value = 108497.75
$ws.Range("L138").Value = 108497.75
$ws.Range("N138").Value = -118777.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 25500302
$ws.Range("I86").Value = 61907664
$ws.Range("J86").Value = 15149.6
$ws.Range("K86").Value = 61907664
$ws.Range("L86").Value = 15149.6
$ws.Range("M86").Value = -61906541
$ws.Range("N86").Value = -17395.6
$ws.Range("H89").Value = 25500302
$ws.Range("I89").Value = 61907664
$ws.Range("J89").Value = 15149.6
$ws.Range("K89").Value = 309538320
$ws.Range("L89").Value = 75748
$ws.Range("M89").Value = -309532704
$ws.Range("N89").Value = -86980
$ws.Range("H105").Value = 2606183.5
$ws.Range("I105").Value = 2719409
$ws.Range("K105").Value = 2719409
$ws.Range("M105").Value = -2717662
$ws.Range("H134").Value = 11418.75
$ws.Range("I134").Value = 9256.27
$ws.Range("K134").Value = 27768.81
$ws.Range("M134").Value = -25233.81

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H122").Value = 1493.8334
$ws.Range("I122").Value = 1568.9412
$ws.Range("J122").Value = 1311.4286
$ws.Range("K122").Value = 4706.8236
$ws.Range("L122").Value = 3934.2858
$ws.Range("M122").Value = -2256.8236
$ws.Range("N122").Value = -8834.2858
$ws.Range("H134").Value = 8905.781000000001
$ws.Range("I134").Value = 7015.4736
$ws.Range("K134").Value = 21046.4208
$ws.Range("M134").Value = -18511.4208
$ws.Range("H138").Value = 29140
$ws.Range("J138").Value = 29140
$ws.Range("L138").Value = 29140
$ws.Range("N138").Value = -39420
$ws.Range("H139").Value = 98354.39999999999
$ws.Range("J139").Value = 98354.39999999999
$ws.Range("L139").Value = 98354.39999999999
$ws.Range("N139").Value = -108634.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 35805.9
$ws.Range("I5").Value = 813.55554
$ws.Range("J5").Value = 93066.09
$ws.Range("K5").Value = 2440.66662
$ws.Range("L5").Value = 279198.27
$ws.Range("M5").Value = -2328.66662
$ws.Range("N5").Value = -279422.27
$ws.Range("H46").Value = 162042.19
$ws.Range("I46").Value = 479150.56
$ws.Range("K46").Value = 1437451.68
$ws.Range("M46").Value = -1437360.68
$ws.Range("H113").Value = 2920
$ws.Range("I113").Value = 4546.846
$ws.Range("K113").Value = 13640.538
$ws.Range("M113").Value = -11470.538
$ws.Range("H132").Value = 1393.6072
$ws.Range("J132").Value = 1253.8422
$ws.Range("L132").Value = 11284.5798
$ws.Range("N132").Value = -16344.5798
$ws.Range("H135").Value = 35805.9
$ws.Range("I135").Value = 813.55554
$ws.Range("J135").Value = 93066.09
$ws.Range("K135").Value = 7321.99986
$ws.Range("L135").Value = 837594.8099999999
$ws.Range("M135").Value = -4786.99986
$ws.Range("N135").Value = -842664.8099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9100020
$ws.Range("I70").Value = 20003802
$ws.Range("J70").Value = 13534.333
$ws.Range("K70").Value = 20003802
$ws.Range("L70").Value = 13534.333
$ws.Range("M70").Value = -20003532
$ws.Range("N70").Value = -14074.333
$ws.Range("H73").Value = 9100020
$ws.Range("I73").Value = 20003802
$ws.Range("J73").Value = 13534.333
$ws.Range("K73").Value = 20003802
$ws.Range("L73").Value = 13534.333
$ws.Range("M73").Value = -20002866
$ws.Range("N73").Value = -15406.333
$ws.Range("H126").Value = 4028058
$ws.Range("I126").Value = 3249136
$ws.Range("K126").Value = 9747408
$ws.Range("M126").Value = -9744938
$ws.Range("H132").Value = 5846.651
$ws.Range("I132").Value = 4537.4253
$ws.Range("J132").Value = 9692.5
$ws.Range("K132").Value = 13612.2759
$ws.Range("L132").Value = 29077.5
$ws.Range("M132").Value = -11082.2759
$ws.Range("N132").Value = -34137.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 4957.75
$ws.Range("I23").Value = 4916.5
$ws.Range("J23").Value = 4999
$ws.Range("K23").Value = 4916.5
$ws.Range("L23").Value = 4999
$ws.Range("M23").Value = -4686.5
$ws.Range("N23").Value = -5459
$ws.Range("H26").Value = 7000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 7000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 7000
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -7590
$ws.Range("H40").Value = 4652.3335
$ws.Range("I40").Value = 2585.2646
$ws.Range("J40").Value = 11041.454
$ws.Range("K40").Value = 2585.2646
$ws.Range("L40").Value = 11041.454
$ws.Range("M40").Value = -2449.2646
$ws.Range("N40").Value = -11313.454
$ws.Range("H61").Value = 22223550
$ws.Range("I61").Value = 22223550
$ws.Range("K61").Value = 22223550
$ws.Range("M61").Value = -22223348
$ws.Range("H93").Value = 7942701
$ws.Range("I93").Value = 12347140
$ws.Range("J93").Value = 14710.4
$ws.Range("K93").Value = 12347140
$ws.Range("L93").Value = 14710.4
$ws.Range("M93").Value = -12345892
$ws.Range("N93").Value = -17206.4
$ws.Range("H113").Value = 22223550
$ws.Range("I113").Value = 22223550
$ws.Range("K113").Value = 22223550
$ws.Range("M113").Value = -22221380
$ws.Range("H122").Value = 5645.952
$ws.Range("I122").Value = 4055.2144
$ws.Range("K122").Value = 12165.6432
$ws.Range("M122").Value = -9715.643199999999
$ws.Range("H132").Value = 12464.529
$ws.Range("I132").Value = 13935.448
$ws.Range("J132").Value = 3933.2
$ws.Range("K132").Value = 41806.344
$ws.Range("L132").Value = 11799.6
$ws.Range("M132").Value = -39276.344
$ws.Range("N132").Value = -16859.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13134925
$ws.Range("I132").Value = 15630957
$ws.Range("J132").Value = 846766.75
$ws.Range("K132").Value = 46892871
$ws.Range("L132").Value = 2540300.25
$ws.Range("M132").Value = -46890341
$ws.Range("N132").Value = -2545360.25
$ws.Range("H136").Value = 5856.316
$ws.Range("I136").Value = 6401.3667
$ws.Range("J136").Value = 3812.375
$ws.Range("K136").Value = 19204.1001
$ws.Range("L136").Value = 11437.125
$ws.Range("M136").Value = -16654.1001
$ws.Range("N136").Value = -16537.125
$ws.Range("H138").Value = 89142.664
$ws.Range("J138").Value = 89142.664
$ws.Range("L138").Value = 89142.664
$ws.Range("N138").Value = -99422.664
$ws.Range("H139").Value = 50374.875
$ws.Range("J139").Value = 50374.875
$ws.Range("L139").Value = 50374.875
$ws.Range("N139").Value = -60654.875
